$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 16 new word-normalization rows (103-118).
# Values are written in this specific order because the target workbook
# assigns new sharedStrings indices in the order cells are set (not row
# order) -- A103 and A113 were entered before the rest, and B113 was
# entered after the main A/B sweep of rows 104-112.
$ws.Range("A103").Value = 'pljrn'
$ws.Range("B103").Value = 'pelajaran'
$ws.Range("A113").Value = 'offline'
$ws.Range("A104").Value = 'pts'
$ws.Range("B104").Value = 'penilaian tengan semester'
$ws.Range("A105").Value = 'mesti'
$ws.Range("B105").Value = 'harus'
$ws.Range("A106").Value = 'emg'
$ws.Range("B106").Value = 'memang'
$ws.Range("A107").Value = 'kalo'
$ws.Range("B107").Value = 'kalau'
$ws.Range("A108").Value = 'cepet'
$ws.Range("B108").Value = 'cepat'
$ws.Range("A109").Value = 'online'
$ws.Range("B109").Value = 'daring'
$ws.Range("A110").Value = 'krn'
$ws.Range("B110").Value = 'karena'
$ws.Range("A111").Value = 'agak'
$ws.Range("B111").Value = 'sedikit'
$ws.Range("A112").Value = 'blm'
$ws.Range("B112").Value = 'belum'
$ws.Range("B113").Value = 'luring'
$ws.Range("A114").Value = 'taun'
$ws.Range("B114").Value = 'tahun'
$ws.Range("A115").Value = 'no'
$ws.Range("B115").Value = 'tidak'
$ws.Range("A116").Value = 'gini'
$ws.Range("B116").Value = 'begini'
$ws.Range("A117").Value = 'full'
$ws.Range("B117").Value = 'penuh'
$ws.Range("A118").Value = 'or'
$ws.Range("B118").Value = 'atau'

# Update selection to the next empty row, matching the final saved view.
$ws.Range("A119").Select()

"Done"
